# "separate dept from affiliations"
#
# PI hours sheet previously had a single "dept" column holding each PI's
# full list of unit affiliations (e.g. "['ME', 'AE', 'CSL']"). The edit:
#   1. Moves that full-affiliation-list column to a new "app" column (F).
#   2. Replaces "dept" (E) with the PI's single primary department.
#   3. Renames "dept hours" -> "department hours" and recomputes it to
#      aggregate only by that single primary department (drops the CSL
#      "unit" bucket, since CSL is cross-cutting, not a primary dept).
#   4. Adds a new sheet "unit(accumulative) hours" holding what used to be
#      the "dept hours" aggregation (every unit each PI is affiliated
#      with, CSL included - i.e. accumulative across units).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "PI hours"
$ws2 = $wb.Worksheets.Item(2)   # "dept hours" -> "department hours"

# --- Step 1: clone the existing "dept hours" sheet (CSL + depts,
# accumulative across all affiliations) to the end of the workbook; this
# becomes the new "unit(accumulative) hours" sheet, untouched otherwise.
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "unit(accumulative) hours"
$ws3.Range("B1").Value = "unit(accumulative)"

# --- Step 2: on "PI hours", introduce the "app" column (F) holding the
# original full affiliation lists, matching the header style used by the
# other header cells.
$ws1.Range("E1").Copy($ws1.Range("F1"))
$ws1.Range("F1").Value = "app"

$ws1.Range("F2").Value = "['ME', 'AE', 'CSL']"
$ws1.Range("F3").Value = "['ECE', 'CSL']"
$ws1.Range("F4").Value = "['ECE', 'CSL']"
$ws1.Range("F5").Value = "['ME', 'CSL']"
$ws1.Range("F6").Value = "['ABE', 'CSL']"
$ws1.Range("F7").Value = "['AE']"

# --- Step 3: "dept" (E) now holds just the PI's single primary department.
$ws1.Range("E2").Value = "ME"
$ws1.Range("E3").Value = "ECE"
$ws1.Range("E4").Value = "ECE"
$ws1.Range("E5").Value = "ME"
$ws1.Range("E6").Value = "ABE"
$ws1.Range("E7").Value = "AE"

# --- Step 4: rename "dept hours" -> "department hours" and recompute it
# so it only aggregates by primary department (CSL row removed).
$ws2.Name = "department hours"

$ws2.Range("B2").Value = "ECE"
$ws2.Range("C2").Value = 52.5
$ws2.Range("D2").Value = 50

$ws2.Range("B3").Value = "ME"
$ws2.Range("C3").Value = 48
$ws2.Range("D3").Value = 45.71428571428572

$ws2.Range("B4").Value = "ABE"
$ws2.Range("C4").Value = 3
$ws2.Range("D4").Value = 2.857142857142857

$ws2.Range("B5").Value = "AE"
$ws2.Range("C5").Value = 1.5
$ws2.Range("D5").Value = 1.428571428571429

$ws2.Rows.Item(6).Delete()

# Restore "PI hours" as the active sheet, as it was originally.
$ws1.Activate()
